$wb = $excel.ActiveWorkbook

$wsNotNorm = $wb.Worksheets.Item("NOT NORMALZIED")
$wsNorm = $wb.Worksheets.Item("NORMALIZED")

# Update raw "Output" values (column D) on the "NOT NORMALZIED" sheet for the
# OR, NAND and NOR gates. E column formulas on that sheet will recalc
# automatically.
$wsNotNorm.Range("D6").Value = 0.473
$wsNotNorm.Range("D7").Value = 3.128
$wsNotNorm.Range("D8").Value = 4.062
$wsNotNorm.Range("D9").Value = 4.188

$wsNotNorm.Range("D14").Value = 3.07
$wsNotNorm.Range("D15").Value = 3.039
$wsNotNorm.Range("D16").Value = 2.899
$wsNotNorm.Range("D17").Value = 0.672

$wsNotNorm.Range("D18").Value = 5.024
$wsNotNorm.Range("D19").Value = 1.135
$wsNotNorm.Range("D20").Value = 1.249
$wsNotNorm.Range("D21").Value = 1.19

$excel.CalculateFull()

# Mirror the recalculated normalized ratios (column E on "NOT NORMALZIED")
# into the static "Output" values (column D) on the "NORMALIZED" sheet.
$wsNorm.Range("D6").Value = $wsNotNorm.Range("E6").Value2
$wsNorm.Range("D7").Value = $wsNotNorm.Range("E7").Value2
$wsNorm.Range("D8").Value = $wsNotNorm.Range("E8").Value2

$wsNorm.Range("D14").Value = $wsNotNorm.Range("E14").Value2
$wsNorm.Range("D15").Value = $wsNotNorm.Range("E15").Value2
$wsNorm.Range("D16").Value = $wsNotNorm.Range("E16").Value2
$wsNorm.Range("D17").Value = $wsNotNorm.Range("E17").Value2

$wsNorm.Range("D19").Value = $wsNotNorm.Range("E19").Value2
$wsNorm.Range("D20").Value = $wsNotNorm.Range("E20").Value2
$wsNorm.Range("D21").Value = $wsNotNorm.Range("E21").Value2

# Update the selection shown on each sheet to match where the author was
# last working (the "NOT NORMALZIED" sheet's OR/NOR rows), leaving the
# "NORMALIZED" sheet as the active tab.
$wsNotNorm.Activate()
$wsNotNorm.Range("E18:E21").Select()

$wsNorm.Activate()
$wsNorm.Range("E18").Select()
